$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.690.95"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "2.221.69"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "240.69"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "0.618"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "74.98"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").Value = "41.28"
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "54.82"
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("D13").Value = "6.89"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "2.554.53"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "14.66"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "2.224.89"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "0.801"
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("D19").Value = "42.519.42"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "70.74"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "5.92"
$ws.Range("E22").Value = "  -4.46%  "
$ws.Range("D23").Value = "10.03"
$ws.Range("E23").Value = "  -9.30%  "
$ws.Range("D24").Value = "229.26"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  +7.10%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("D28").Value = "3.41"
$ws.Range("E28").Value = "  -6.02%  "
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "172.59"
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.09"
$ws.Range("E31").Value = "  -4.85%  "
$ws.Range("D32").Value = "35.82"
$ws.Range("E32").Value = "  +16.69%  "
$ws.Range("D33").Value = "20.30"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").Value = "0.0794"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.40"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.107"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").Value = "0.0322"
$ws.Range("E39").Value = "  +6.25%  "
$ws.Range("D40").Value = "12.49"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").Value = "2.14"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "5.51"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("D43").Value = "60.78"
$ws.Range("E43").Value = "  -5.57%  "
$ws.Range("D44").Value = "0.198"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "8.56"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "0.0986"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").Value = "99.32"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "2.30"
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  -3.14%  "
$ws.Range("D51").Value = "0.420"
$ws.Range("E51").Value = "  +13.49%  "
